$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1298.6
$ws.Range("J40").Value = 1298.6
$ws.Range("L40").Value = 1298.6
$ws.Range("N40").Value = -1648.6

$ws.Range("H96").Value = 463.8889
$ws.Range("I96").Value = 171.875
$ws.Range("K96").Value = 515.625
$ws.Range("M96").Value = 857.375

$ws.Range("H99").Value = 614.2857
$ws.Range("I99").Value = 383.33334
$ws.Range("K99").Value = 1150.00002
$ws.Range("M99").Value = 347.9999800000001

$ws.Range("H132").Value = 857.25
$ws.Range("I132").Value = 857.25
$ws.Range("K132").Value = 2571.75
$ws.Range("M132").Value = -41.75

$ws.Range("H137").Value = 10002
$ws.Range("I137").Value = 10002
$ws.Range("K137").Value = 30006
$ws.Range("M137").Value = -27456

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 5025
$ws.Range("I34").Value = 5025
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 5025
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -4754
$ws.Range("N34").ClearContents()

$ws.Range("H38").Value = 1247617.2
$ws.Range("I38").Value = 3234.5
$ws.Range("J38").Value = 2492000
$ws.Range("K38").Value = 3234.5
$ws.Range("L38").Value = 2492000
$ws.Range("M38").Value = -2767.5
$ws.Range("N38").Value = -2492934

$ws.Range("H63").Value = 928.5
$ws.Range("I63").Value = 928.5
$ws.Range("K63").Value = 928.5
$ws.Range("M63").Value = -242.5

$ws.Range("H66").Value = 928.5
$ws.Range("I66").Value = 928.5
$ws.Range("K66").Value = 4642.5
$ws.Range("M66").Value = -1210.5

$ws.Range("H124").Value = 82000
$ws.Range("J124").Value = 82000
$ws.Range("L124").Value = 82000
$ws.Range("N124").Value = -91820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H75").Value = 5000
$ws.Range("I75").Value = 5000
$ws.Range("K75").Value = 5000
$ws.Range("M75").Value = -4064

$ws.Range("H78").Value = 5000
$ws.Range("I78").Value = 5000
$ws.Range("K78").Value = 15000
$ws.Range("M78").Value = -10320

$ws.Range("H134").Value = 1563.4286
$ws.Range("I134").Value = 1657.3334
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 4972.0002
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -2437.0002
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 399.42856
$ws.Range("I22").Value = 416
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 416
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -66
$ws.Range("N22").Value = -1000

$ws.Range("H33").Value = 41941.5
$ws.Range("I33").Value = 11836.2
$ws.Range("K33").Value = 11836.2
$ws.Range("M33").Value = -11457.2

$ws.Range("H62").Value = 4099.25
$ws.Range("I62").Value = 4099.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4099.25
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3475.25
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 4099.25
$ws.Range("I65").Value = 4099.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 20496.25
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -17376.25
$ws.Range("N65").ClearContents()

$ws.Range("H134").Value = 4417.5
$ws.Range("I134").Value = 1722
$ws.Range("J134").Value = 11156.25
$ws.Range("K134").Value = 5166
$ws.Range("L134").Value = 33468.75
$ws.Range("M134").Value = -2631
$ws.Range("N134").Value = -38538.75

$ws.Range("H141").Value = 122224
$ws.Range("J141").Value = 122224
$ws.Range("L141").Value = 122224
$ws.Range("N141").Value = -132584

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2123.75
$ws.Range("I51").Value = 1498.3334
$ws.Range("K51").Value = 4495.0002
$ws.Range("M51").Value = -4035.0002

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3965.6667
$ws.Range("I80").Value = 3898
$ws.Range("K80").Value = 3898
$ws.Range("M80").Value = -2900

$ws.Range("H83").Value = 3965.6667
$ws.Range("I83").Value = 3898
$ws.Range("K83").Value = 19490
$ws.Range("M83").Value = -14498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 10503.75
$ws.Range("J31").Value = 16000
$ws.Range("L31").Value = 16000
$ws.Range("N31").Value = -16496

$ws.Range("H46").Value = 3212
$ws.Range("I46").Value = 949.3333
$ws.Range("K46").Value = 949.3333
$ws.Range("M46").Value = -761.3333

$ws.Range("H68").Value = 1985.7142
$ws.Range("I68").Value = 1985.7142
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1985.7142
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1236.7142
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1985.7142
$ws.Range("I71").Value = 1985.7142
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9928.571
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6184.571
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2636.8
$ws.Range("I81").Value = 2981.1428
$ws.Range("J81").Value = 1833.3334
$ws.Range("K81").Value = 5962.2856
$ws.Range("L81").Value = 3666.6668
$ws.Range("M81").Value = -4901.2856
$ws.Range("N81").Value = -5788.6668

$ws.Range("H84").Value = 2636.8
$ws.Range("I84").Value = 2981.1428
$ws.Range("J84").Value = 1833.3334
$ws.Range("K84").Value = 29811.428
$ws.Range("L84").Value = 18333.334
$ws.Range("M84").Value = -24507.428
$ws.Range("N84").Value = -28941.334

$ws.Range("H132").Value = 944.375
$ws.Range("I132").Value = 936.4286
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 2809.2858
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -279.2857999999997
$ws.Range("N132").Value = -8060
